# Auto-generated edit script applying numeric corrections to the
# Seraph_Profits leve-crafting profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4004.2
$ws.Range("I86").Value = 2366.5
$ws.Range("K86").Value = 2366.5
$ws.Range("M86").Value = -1243.5
$ws.Range("H89").Value = 4004.2
$ws.Range("I89").Value = 2366.5
$ws.Range("K89").Value = 11832.5
$ws.Range("M89").Value = -6216.5
$ws.Range("H141").Value = 5296.5557
$ws.Range("I141").Value = 5278.3335
$ws.Range("K141").Value = 15835.0005
$ws.Range("M141").Value = -10655.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 174008.17
$ws.Range("I6").Value = 999999
$ws.Range("J6").Value = 8810
$ws.Range("K6").Value = 999999
$ws.Range("L6").Value = 8810
$ws.Range("M6").Value = -999826
$ws.Range("N6").Value = -9156
$ws.Range("H30").Value = 375
$ws.Range("I30").Value = 90
$ws.Range("K30").Value = 90
$ws.Range("M30").Value = 60
$ws.Range("H32").Value = 14720.937
$ws.Range("I32").Value = 6337.054
$ws.Range("J32").Value = 26651.846
$ws.Range("K32").Value = 6337.054
$ws.Range("L32").Value = 26651.846
$ws.Range("M32").Value = -6050.054
$ws.Range("N32").Value = -27225.846
$ws.Range("H61").Value = 1498.5
$ws.Range("I61").Value = 1498.5
$ws.Range("K61").Value = 1498.5
$ws.Range("M61").Value = -1286.5
$ws.Range("H74").Value = 5191.375
$ws.Range("I74").Value = 1598.75
$ws.Range("K74").Value = 1598.75
$ws.Range("M74").Value = -724.75
$ws.Range("H77").Value = 5191.375
$ws.Range("I77").Value = 1598.75
$ws.Range("K77").Value = 7993.75
$ws.Range("M77").Value = -3625.75
$ws.Range("H102").Value = 1570.9
$ws.Range("I102").Value = 958.4286
$ws.Range("K102").Value = 958.4286
$ws.Range("M102").Value = 663.5714
$ws.Range("H136").Value = 1498.5
$ws.Range("I136").Value = 1498.5
$ws.Range("K136").Value = 4495.5
$ws.Range("M136").Value = -1945.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 742
$ws.Range("J16").Value = 1499.5
$ws.Range("L16").Value = 1499.5
$ws.Range("N16").Value = -2073.5
$ws.Range("H32").Value = 500750
$ws.Range("I32").Value = 500750
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 500750
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -500434
$ws.Range("N32").ClearContents()
$ws.Range("H37").Value = 5475
$ws.Range("I37").Value = 5400
$ws.Range("J37").Value = 5500
$ws.Range("K37").Value = 5400
$ws.Range("L37").Value = 5500
$ws.Range("M37").Value = -5293
$ws.Range("N37").Value = -5714
$ws.Range("H107").Value = 418.55554
$ws.Range("I107").Value = 314.625
$ws.Range("K107").Value = 314.625
$ws.Range("M107").Value = 1605.375
$ws.Range("H113").Value = 742
$ws.Range("J113").Value = 1499.5
$ws.Range("L113").Value = 1499.5
$ws.Range("N113").Value = -5839.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 686.25
$ws.Range("I5").Value = 413.85715
$ws.Range("J5").Value = 1067.6
$ws.Range("K5").Value = 1241.57145
$ws.Range("L5").Value = 3202.8
$ws.Range("M5").Value = -1129.57145
$ws.Range("N5").Value = -3426.8
$ws.Range("H113").Value = 848.7353000000001
$ws.Range("J113").Value = 748.7143
$ws.Range("L113").Value = 2246.1429
$ws.Range("N113").Value = -6586.1429
$ws.Range("H135").Value = 686.25
$ws.Range("I135").Value = 413.85715
$ws.Range("J135").Value = 1067.6
$ws.Range("K135").Value = 3724.71435
$ws.Range("L135").Value = 9608.4
$ws.Range("M135").Value = -1189.71435
$ws.Range("N135").Value = -14678.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18000
$ws.Range("J15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("N15").Value = -18576
$ws.Range("H80").Value = 4373.25
$ws.Range("I80").Value = 3286.25
$ws.Range("K80").Value = 3286.25
$ws.Range("M80").Value = -2288.25
$ws.Range("H81").Value = 18000
$ws.Range("J81").Value = 18000
$ws.Range("L81").Value = 18000
$ws.Range("N81").Value = -19996
$ws.Range("H83").Value = 4373.25
$ws.Range("I83").Value = 3286.25
$ws.Range("K83").Value = 16431.25
$ws.Range("M83").Value = -11439.25
$ws.Range("H84").Value = 18000
$ws.Range("J84").Value = 18000
$ws.Range("L84").Value = 54000
$ws.Range("N84").Value = -63984
$ws.Range("H126").Value = 5376.375
$ws.Range("I126").Value = 3012
$ws.Range("J126").Value = 5714.143
$ws.Range("K126").Value = 9036
$ws.Range("L126").Value = 17142.429
$ws.Range("M126").Value = -6566
$ws.Range("N126").Value = -22082.429
$ws.Range("H132").Value = 3376.5
$ws.Range("I132").Value = 3035.7222
$ws.Range("J132").Value = 6443.5
$ws.Range("K132").Value = 9107.1666
$ws.Range("L132").Value = 19330.5
$ws.Range("M132").Value = -6577.1666
$ws.Range("N132").Value = -24390.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2439.2
$ws.Range("I7").Value = 2439.2
$ws.Range("K7").Value = 2439.2
$ws.Range("M7").Value = -2327.2
$ws.Range("H16").Value = 4256.85
$ws.Range("J16").Value = 5863.3335
$ws.Range("L16").Value = 5863.3335
$ws.Range("N16").Value = -6203.3335
$ws.Range("H46").Value = 2731.0527
$ws.Range("I46").Value = 1988.3334
$ws.Range("J46").Value = 3399.5
$ws.Range("K46").Value = 1988.3334
$ws.Range("L46").Value = 3399.5
$ws.Range("M46").Value = -1800.3334
$ws.Range("N46").Value = -3775.5
$ws.Range("H68").Value = 2997.6667
$ws.Range("I68").Value = 2994
$ws.Range("K68").Value = 2994
$ws.Range("M68").Value = -2245
$ws.Range("H71").Value = 2997.6667
$ws.Range("I71").Value = 2994
$ws.Range("K71").Value = 14970
$ws.Range("M71").Value = -11226
$ws.Range("H126").Value = 2439.2
$ws.Range("I126").Value = 2439.2
$ws.Range("K126").Value = 7317.599999999999
$ws.Range("M126").Value = -4847.599999999999
$ws.Range("H132").Value = 5197.6484
$ws.Range("I132").Value = 4828.3335
$ws.Range("J132").Value = 5879.4614
$ws.Range("K132").Value = 14485.0005
$ws.Range("L132").Value = 17638.3842
$ws.Range("M132").Value = -11955.0005
$ws.Range("N132").Value = -22698.3842
$ws.Range("H136").Value = 3445.6316
$ws.Range("I136").Value = 3310.5625
$ws.Range("K136").Value = 9931.6875
$ws.Range("M136").Value = -7381.6875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 936.4286
$ws.Range("I107").Value = 425.83334
$ws.Range("K107").Value = 1277.50002
$ws.Range("M107").Value = 642.4999800000001
$ws.Range("H122").Value = 1292.9333
$ws.Range("I122").Value = 1292.9333
$ws.Range("K122").Value = 3878.7999
$ws.Range("M122").Value = -1428.7999
$ws.Range("H123").Value = 48999.5
$ws.Range("J123").Value = 48999.5
$ws.Range("L123").Value = 48999.5
$ws.Range("N123").Value = -58799.5
$ws.Range("H126").Value = 73884.92999999999
$ws.Range("I126").Value = 84782.414
$ws.Range("K126").Value = 254347.242
$ws.Range("M126").Value = -251877.242
$ws.Range("H132").Value = 1567.6666
$ws.Range("I132").Value = 1562.7778
$ws.Range("K132").Value = 4688.3334
$ws.Range("M132").Value = -2158.3334
$ws.Range("H136").Value = 64419.5
$ws.Range("I136").Value = 1250.909
$ws.Range("K136").Value = 3752.727
$ws.Range("M136").Value = -1202.727

